$wb = $excel.ActiveWorkbook

# ALC row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1727.0667
$ws.Range("I53").Value = 443
$ws.Range("J53").Value = 2194
$ws.Range("K53").Value = 443
$ws.Range("L53").Value = 2194
$ws.Range("M53").Value = 194
$ws.Range("N53").Value = -3468

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 15774.777
$ws.Range("I76").Value = 14852.286
$ws.Range("J76").Value = 19003.5
$ws.Range("K76").Value = 14852.286
$ws.Range("L76").Value = 19003.5
$ws.Range("M76").Value = -14537.286
$ws.Range("N76").Value = -19633.5

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 15774.777
$ws.Range("I79").Value = 14852.286
$ws.Range("J79").Value = 19003.5
$ws.Range("K79").Value = 14852.286
$ws.Range("L79").Value = 19003.5
$ws.Range("M79").Value = -13760.286
$ws.Range("N79").Value = -21187.5

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1823.6666
$ws.Range("I80").Value = 491.86667
$ws.Range("J80").Value = 4043.3333
$ws.Range("K80").Value = 1475.60001
$ws.Range("L80").Value = 12129.9999
$ws.Range("M80").Value = -477.6000100000001
$ws.Range("N80").Value = -14125.9999

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1823.6666
$ws.Range("I83").Value = 491.86667
$ws.Range("J83").Value = 4043.3333
$ws.Range("K83").Value = 4426.80003
$ws.Range("L83").Value = 36389.9997
$ws.Range("M83").Value = 565.1999699999997
$ws.Range("N83").Value = -46373.9997

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 13993.818
$ws.Range("I88").Value = 3500
$ws.Range("J88").Value = 16325.777
$ws.Range("K88").Value = 3500
$ws.Range("L88").Value = 16325.777
$ws.Range("M88").Value = -3094
$ws.Range("N88").Value = -17137.777

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 13993.818
$ws.Range("I91").Value = 3500
$ws.Range("J91").Value = 16325.777
$ws.Range("K91").Value = 3500
$ws.Range("L91").Value = 16325.777
$ws.Range("M91").Value = -2096
$ws.Range("N91").Value = -19133.777

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3514.5
$ws.Range("J100").Value = 3633.5833
$ws.Range("L100").Value = 3633.5833
$ws.Range("N100").Value = -4715.5833

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1848.9166
$ws.Range("J103").Value = 2097.889
$ws.Range("L103").Value = 6293.667
$ws.Range("N103").Value = -7465.667

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3460.2307
$ws.Range("I112").Value = 373.625
$ws.Range("J112").Value = 8398.799999999999
$ws.Range("K112").Value = 1120.875
$ws.Range("L112").Value = 25196.4
$ws.Range("M112").Value = -12.875
$ws.Range("N112").Value = -27412.4

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 957.1429000000001
$ws.Range("I125").Value = 946.8333
$ws.Range("J125").Value = 964.875
$ws.Range("K125").Value = 8521.4997
$ws.Range("L125").Value = 8683.875
$ws.Range("M125").Value = -6061.4997
$ws.Range("N125").Value = -13603.875

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1074.0975
$ws.Range("I132").Value = 991.10254
$ws.Range("J132").Value = 2692.5
$ws.Range("K132").Value = 2973.30762
$ws.Range("L132").Value = 8077.5
$ws.Range("M132").Value = -443.30762
$ws.Range("N132").Value = -13137.5

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2756.1738
$ws.Range("I137").Value = 851
$ws.Range("K137").Value = 2553
$ws.Range("M137").Value = -3

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2657.7659
$ws.Range("J138").Value = 3164.7307
$ws.Range("L138").Value = 9494.1921
$ws.Range("N138").Value = -19774.1921

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3577.7646
$ws.Range("I132").Value = 2841.1667
$ws.Range("J132").Value = 5345.6
$ws.Range("K132").Value = 8523.500100000001
$ws.Range("L132").Value = 16036.8
$ws.Range("M132").Value = -5993.500100000001
$ws.Range("N132").Value = -21096.8

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2950.6667
$ws.Range("I20").Value = 2224.3
$ws.Range("J20").Value = 4403.4
$ws.Range("K20").Value = 2224.3
$ws.Range("L20").Value = 4403.4
$ws.Range("M20").Value = -1977.3
$ws.Range("N20").Value = -4897.4

# BSM row 28
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 39978
$ws.Range("J28").Value = 39978
$ws.Range("L28").Value = 39978
$ws.Range("N28").Value = -40566

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5518.5933
$ws.Range("I132").Value = 5325.8335
$ws.Range("J132").Value = 7600.4
$ws.Range("K132").Value = 15977.5005
$ws.Range("L132").Value = 22801.2
$ws.Range("M132").Value = -13447.5005
$ws.Range("N132").Value = -27861.2

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 33483.934
$ws.Range("I7").Value = 121.416664
$ws.Range("K7").Value = 364.249992
$ws.Range("M7").Value = -252.249992

# CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 344547.75
$ws.Range("I9").Value = 577134.4
$ws.Range("J9").Value = 4613.4614
$ws.Range("K9").Value = 1731403.2
$ws.Range("L9").Value = 13840.3842
$ws.Range("M9").Value = -1731179.2
$ws.Range("N9").Value = -14288.3842

# CUL row 47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 17131.666
$ws.Range("I47").Value = 17131.666
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 51394.99800000001
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -50963.99800000001
$ws.Range("N47").ClearContents()

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13413.148
$ws.Range("I70").Value = 5461.2104
$ws.Range("J70").Value = 32299
$ws.Range("K70").Value = 5461.2104
$ws.Range("L70").Value = 32299
$ws.Range("M70").Value = -5191.2104
$ws.Range("N70").Value = -32839

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13413.148
$ws.Range("I73").Value = 5461.2104
$ws.Range("J73").Value = 32299
$ws.Range("K73").Value = 5461.2104
$ws.Range("L73").Value = 32299
$ws.Range("M73").Value = -4525.2104
$ws.Range("N73").Value = -34171

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2733.48
$ws.Range("I122").Value = 2397.7222
$ws.Range("K122").Value = 7193.1666
$ws.Range("M122").Value = -4743.1666

# GSM row 125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920

# LTW row 6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 63995
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 63995
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 63995
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -64219

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2462.5881
$ws.Range("I16").Value = 2190.9333
$ws.Range("K16").Value = 2190.9333
$ws.Range("M16").Value = -2020.9333

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8434.210999999999
$ws.Range("I40").Value = 8316
$ws.Range("J40").Value = 8877.5
$ws.Range("K40").Value = 8316
$ws.Range("L40").Value = 8877.5
$ws.Range("M40").Value = -8180
$ws.Range("N40").Value = -9149.5

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2315.2354
$ws.Range("I93").Value = 1939.8
$ws.Range("J93").Value = 3358.111
$ws.Range("K93").Value = 1939.8
$ws.Range("L93").Value = 3358.111
$ws.Range("M93").Value = -691.8
$ws.Range("N93").Value = -5854.111

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7343.778
$ws.Range("I132").Value = 1949.75
$ws.Range("J132").Value = 11659
$ws.Range("K132").Value = 5849.25
$ws.Range("L132").Value = 34977
$ws.Range("M132").Value = -3319.25
$ws.Range("N132").Value = -40037

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6996.6
$ws.Range("I136").Value = 4170.375
$ws.Range("K136").Value = 12511.125
$ws.Range("M136").Value = -9961.125

# WVR row 99
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 60000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 503
$ws.Range("I107").Value = 488.69232
$ws.Range("J107").Value = 534
$ws.Range("K107").Value = 1466.07696
$ws.Range("L107").Value = 1602
$ws.Range("M107").Value = 453.9230400000001
$ws.Range("N107").Value = -5442

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4532.231
$ws.Range("I132").Value = 4183.375
$ws.Range("J132").Value = 6127
$ws.Range("K132").Value = 12550.125
$ws.Range("L132").Value = 18381
$ws.Range("M132").Value = -10020.125
$ws.Range("N132").Value = -23441

# WVR row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 69997.2
$ws.Range("J137").Value = 69997.2
$ws.Range("L137").Value = 69997.2
$ws.Range("N137").Value = -80197.2
